$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the columns that were removed from the report:
# D:F = 21_FL_TTCA, 22_FL_umami, 23_FL_AMAP_alternative
# I   = 26_FL_MRSA (after D:F removal this becomes column F, but we delete by
#       original letter reference before the shift takes effect for the next call)
$ws.Range("D1:F1").EntireColumn.Delete()
$ws.Range("F1").EntireColumn.Delete()
